# Weekly update: insert a new price record for Feria Lagunitas de Puerto
# Montt - Espárragos as the new row 20 (pushing the existing rows 20-44
# down to 21-45), matching the "Fruta / hortaliza, semanal" refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing data rows down by inserting a fresh row at position 20.
$ws.Rows(20).Insert()

# Populate the newly inserted row with this week's record.
$ws.Range("A20").Value = 4
$ws.Range("B20").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C20").Value = "Los Lagos"
$ws.Range("D20").Value = 44540
$ws.Range("E20").Value = 10
$ws.Range("F20").Value = 300000000
$ws.Range("G20").Value = "Espárragos"
$ws.Range("H20").Value = "Sin especificar"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 600
$ws.Range("K20").Value = 1700
$ws.Range("L20").Value = 1700
$ws.Range("M20").Value = 1700
$ws.Range("N20").Value = "`$/kilo"
$ws.Range("O20").Value = "Provincia de Linares"
$ws.Range("P20").Value = 1700
$ws.Range("Q20").Value = 1
$ws.Range("R20").Value = "Hortaliza"
